$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-06 Tuesday" "2024-08-07 Wednesday"

Replace-Text "714÷6=" "485÷7="
Replace-Text "137÷9=" "243÷2="
Replace-Text "666÷2=" "558÷5="
Replace-Text "712÷9=" "115÷4="
Replace-Text "661÷3=" "532÷6="
Replace-Text "675÷5=" "819÷2="
Replace-Text "205÷4=" "361÷3="
Replace-Text "806÷5=" "328÷7="
Replace-Text "528÷4=" "632÷3="
Replace-Text "258÷3=" "160÷6="
Replace-Text "605÷5=" "679÷9="
Replace-Text "352÷3=" "151÷2="
Replace-Text "459÷9=" "526÷4="
Replace-Text "880÷7=" "463÷6="
Replace-Text "774÷4=" "687÷3="
Replace-Text "940÷7=" "400÷3="
Replace-Text "325÷7=" "189÷8="
Replace-Text "962÷2=" "679÷8="
Replace-Text "234÷9=" "294÷2="
Replace-Text "604÷6=" "374÷4="
Replace-Text "640÷6=" "414÷6="
Replace-Text "845÷7=" "814÷7="
Replace-Text "849÷4=" "825÷8="
Replace-Text "363÷8=" "980÷6="
Replace-Text "332÷6=" "149÷3="
